$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 900
$ws.Range("I31").Value = 900
$ws.Range("K31").Value = 2700
$ws.Range("M31").Value = -2470
$ws.Range("H41").Value = 287.26666
$ws.Range("I41").Value = 145.125
$ws.Range("J41").Value = 449.7143
$ws.Range("K41").Value = 145.125
$ws.Range("L41").Value = 449.7143
$ws.Range("M41").Value = 294.875
$ws.Range("N41").Value = -1329.7143
$ws.Range("H94").Value = 7993.3335
$ws.Range("I94").Value = 7993.3335
$ws.Range("K94").Value = 7993.3335
$ws.Range("M94").Value = -7542.3335
$ws.Range("H103").Value = 1010912.8
$ws.Range("I103").Value = 2222748.2
$ws.Range("J103").Value = 1050
$ws.Range("K103").Value = 6668244.600000001
$ws.Range("L103").Value = 3150
$ws.Range("M103").Value = -6667658.600000001
$ws.Range("N103").Value = -4322

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 716432.3
$ws.Range("I33").Value = 716432.3
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 716432.3
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -716103.3
$ws.Range("N33").Value = $null
$ws.Range("H88").Value = 2339.5334
$ws.Range("I88").Value = 1834.3334
$ws.Range("J88").Value = 2676.3333
$ws.Range("K88").Value = 1834.3334
$ws.Range("L88").Value = 2676.3333
$ws.Range("M88").Value = -1428.3334
$ws.Range("N88").Value = -3488.3333
$ws.Range("H91").Value = 2339.5334
$ws.Range("I91").Value = 1834.3334
$ws.Range("J91").Value = 2676.3333
$ws.Range("K91").Value = 1834.3334
$ws.Range("L91").Value = 2676.3333
$ws.Range("M91").Value = -430.3334
$ws.Range("N91").Value = -5484.3333
$ws.Range("H102").Value = 1701.4546
$ws.Range("I102").Value = 1563.6842
$ws.Range("J102").Value = 2574
$ws.Range("K102").Value = 1563.6842
$ws.Range("L102").Value = 2574
$ws.Range("M102").Value = 58.31580000000008
$ws.Range("N102").Value = -5818
$ws.Range("H110").Value = 1588.0435
$ws.Range("I110").Value = 685
$ws.Range("J110").Value = 2282.6924
$ws.Range("K110").Value = 685
$ws.Range("L110").Value = 2282.6924
$ws.Range("M110").Value = 1360
$ws.Range("N110").Value = -6372.6924
$ws.Range("H122").Value = 2024.6
$ws.Range("I122").Value = 1646.1765
$ws.Range("K122").Value = 4938.529500000001
$ws.Range("M122").Value = -2488.529500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1450
$ws.Range("I86").Value = 1387.5
$ws.Range("J86").Value = 1575
$ws.Range("K86").Value = 1387.5
$ws.Range("L86").Value = 1575
$ws.Range("M86").Value = -264.5
$ws.Range("N86").Value = -3821
$ws.Range("H89").Value = 1450
$ws.Range("I89").Value = 1387.5
$ws.Range("J89").Value = 1575
$ws.Range("K89").Value = 6937.5
$ws.Range("L89").Value = 7875
$ws.Range("M89").Value = -1321.5
$ws.Range("N89").Value = -19107

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = $null
$ws.Range("N33").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 634.3333
$ws.Range("I98").Value = 551.5
$ws.Range("J98").Value = 800
$ws.Range("K98").Value = 1654.5
$ws.Range("L98").Value = 2400
$ws.Range("M98").Value = -156.5
$ws.Range("N98").Value = -5396

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1118.7241
$ws.Range("I97").Value = 935.45
$ws.Range("J97").Value = 1526
$ws.Range("K97").Value = 935.45
$ws.Range("L97").Value = 1526
$ws.Range("M97").Value = -439.45
$ws.Range("N97").Value = -2518
$ws.Range("H132").Value = 5228.8125
$ws.Range("I132").Value = 5705.2085
$ws.Range("K132").Value = 17115.6255
$ws.Range("M132").Value = -14585.6255

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2500.5715
$ws.Range("I7").Value = 1901.3334
$ws.Range("J7").Value = 2950
$ws.Range("K7").Value = 1901.3334
$ws.Range("L7").Value = 2950
$ws.Range("M7").Value = -1789.3334
$ws.Range("N7").Value = -3174
$ws.Range("H30").Value = 691.6667
$ws.Range("I30").Value = 691.6667
$ws.Range("K30").Value = 691.6667
$ws.Range("M30").Value = -583.6667
$ws.Range("H40").Value = 1837.3549
$ws.Range("I40").Value = 1885.1305
$ws.Range("J40").Value = 1700
$ws.Range("K40").Value = 1885.1305
$ws.Range("L40").Value = 1700
$ws.Range("M40").Value = -1749.1305
$ws.Range("N40").Value = -1972
$ws.Range("H100").Value = 3158.5833
$ws.Range("J100").Value = 3375
$ws.Range("L100").Value = 3375
$ws.Range("N100").Value = -4457
$ws.Range("H122").Value = 2776.2856
$ws.Range("I122").Value = 2566.8
$ws.Range("J122").Value = 3300
$ws.Range("K122").Value = 7700.400000000001
$ws.Range("L122").Value = 9900
$ws.Range("M122").Value = -5250.400000000001
$ws.Range("N122").Value = -14800
$ws.Range("H126").Value = 2500.5715
$ws.Range("I126").Value = 1901.3334
$ws.Range("J126").Value = 2950
$ws.Range("K126").Value = 5704.0002
$ws.Range("L126").Value = 8850
$ws.Range("M126").Value = -3234.0002
$ws.Range("N126").Value = -13790

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 14600
$ws.Range("J28").Value = 14600
$ws.Range("L28").Value = 14600
$ws.Range("N28").Value = -15296
$ws.Range("H63").Value = 27394.111
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 27394.111
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 27394.111
$ws.Range("M63").Value = $null
$ws.Range("N63").Value = -28642.111
$ws.Range("H66").Value = 27394.111
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 27394.111
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 82182.333
$ws.Range("M66").Value = $null
$ws.Range("N66").Value = -88422.333
$ws.Range("H81").Value = 4819.9707
$ws.Range("I81").Value = 7686.875
$ws.Range("J81").Value = 2271.611
$ws.Range("K81").Value = 15373.75
$ws.Range("L81").Value = 4543.222
$ws.Range("M81").Value = -14312.75
$ws.Range("N81").Value = -6665.222
$ws.Range("H84").Value = 4819.9707
$ws.Range("I84").Value = 7686.875
$ws.Range("J84").Value = 2271.611
$ws.Range("K84").Value = 76868.75
$ws.Range("L84").Value = 22716.11
$ws.Range("M84").Value = -71564.75
$ws.Range("N84").Value = -33324.11
$ws.Range("H100").Value = 7577168.5
$ws.Range("I100").Value = 15152408
$ws.Range("J100").Value = 1928.8334
$ws.Range("K100").Value = 30304816
$ws.Range("L100").Value = 3857.6668
$ws.Range("M100").Value = -30304275
$ws.Range("N100").Value = -4939.6668
